$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: remove the 张悦 row (row 2, all-zero row) and shift remaining rows up ---
# Current rows 3,4,5 (卢楠, 冷雪, 屈昂) become rows 2,3,4.
$ws1.Rows.Item(2).Delete()

# Update group name column (A) from 质控组 to 北京组 for the 3 remaining data rows
$ws1.Range("A2:A4").Value = "北京组"

# --- Sheet2: update group name column (A) from 质控组 to 北京组 ---
$ws2.Range("A2").Value = "北京组"

# --- Make Sheet2 the active/selected sheet, matching tabSelected move ---
[void]$ws1.Range("A4").Select()
[void]$ws2.Activate()
[void]$ws2.Range("A3").Select()
